# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 691
$ws1.Range("F3").Value = 31
$ws1.Range("F4").Value = 232
$ws1.Range("F5").Value = 2117
$ws1.Range("F6").Value = 46
$ws1.Range("F7").Value = 3416
$ws1.Range("F9").Value = 845

# --- Sheet "全部类型" (fourth sheet) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 691
$ws4.Range("F3").Value = 31
$ws4.Range("F5").Value = 232
$ws4.Range("F6").Value = 2117
$ws4.Range("F7").Value = 46
$ws4.Range("F8").Value = 3416
$ws4.Range("F10").Value = 845
